# Apply the "Added two UART opcodes to control the Arduino IO pins" edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("XMOS->STM32")

# New command names first (so shared-string indices match the source order)
$ws.Range("A37").Value = "Get IO pin status"
$ws.Range("A38").Value = "Set IO pin state"

# Then the new data-width descriptions
$ws.Range("D37").Value = "5-bit cmd, 17-bit reply"
$ws.Range("D38").Value = "17-bit"

# Update the active/selected cell shown when the sheet is opened
$ws.Activate()
$ws.Range("H33").Select()
